# Add new Borrelia burgdorferi B31 control datasets at 5, 10, and 25 reads.
# Relates to issue #100.
#
# The three new records are inserted as rows 5-7 (pushing the previously
# existing rows 5-9 down to rows 8-12), and the CSID lettering sequence in
# column D is extended to keep every record's short id unique.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 5; formatting is inherited from
# the row above, matching the rest of the "Borrelia_burgdorferi_B31_13" block.
$ws.Rows("5:7").Insert()

# Columns that are identical across the new rows and already reuse
# existing values elsewhere in the sheet.
$ws.Range("A5").Value = "Borrelia_burgdorferi_B31_13"
$ws.Range("A6").Value = "Borrelia_burgdorferi_B31_13"
$ws.Range("A7").Value = "Borrelia_burgdorferi_B31_13"
$ws.Range("E5").Value = "Unknown"
$ws.Range("E6").Value = "Unknown"
$ws.Range("E7").Value = "Unknown"
$ws.Range("F5").Value = "Unknown"
$ws.Range("F6").Value = "Unknown"
$ws.Range("F7").Value = "Unknown"
$ws.Range("G5").Value = "NA"
$ws.Range("G6").Value = "NA"
$ws.Range("G7").Value = "NA"
$ws.Range("H5").Value = "NA"
$ws.Range("H6").Value = "NA"
$ws.Range("H7").Value = "NA"
$ws.Range("J5").Value = "a"
$ws.Range("J6").Value = "a"
$ws.Range("J7").Value = "a"

# Row 7: 5-reads control.
$ws.Range("C7").Value = "Control_Borrelia_burgdoreri_B31_5"
# Row 6: 10-reads control.
$ws.Range("D6").Value = "i"
$ws.Range("I7").Value = "Borrelia_burgdoreri_B31_5"
$ws.Range("C6").Value = "Control_Borrelia_burgdoreri_B31_10"
$ws.Range("I6").Value = "Borrelia_burgdoreri_B31_10"
$ws.Range("D7").Value = "j"

# The rows that used to be 5-9 are now 8-12; continue the CSID lettering
# sequence (k, l) for them.
$ws.Range("D8").Value = "k"
$ws.Range("D9").Value = "l"

# Row 5: 25-reads control.
$ws.Range("C5").Value = "Control_Borrelia_burgdoreri_B31_25"
$ws.Range("I5").Value = "Borrelia_burgdoreri_B31_25"

# Continue the CSID lettering sequence (m, n, o) for the rest of the
# shifted rows.
$ws.Range("D10").Value = "m"
$ws.Range("D11").Value = "n"
$ws.Range("D12").Value = "o"

# Row 5's CSID reuses the letter "f" that was freed up by this reshuffle.
$ws.Range("D5").Value = "f"

[void]$ws.Range("D19").Select()
